$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ramírez Buendía Daniel Sebastián "
$ws.Range("B2").Value = 165843
$ws.Range("C2").Value = "1714108568"
$ws.Range("D2").Value = "Ramírez Buendía Daniel Sebastián "
$ws.Range("E2").Value = "156341"
$ws.Range("F2").Value = Get-Date -Year 2020 -Month 11 -Day 6
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Value = "$ 200"

$ws.Range("A3").Value = "Chasiloa Páez Mirian Amparo"
$ws.Range("B3").Value = 45872
$ws.Range("C3").Value = "1003834627001"
$ws.Range("D3").Value = " Yithos"
$ws.Range("E3").Value = "45872"
$ws.Range("F3").Value = Get-Date -Year 2020 -Month 11 -Day 9
$ws.Range("F3").NumberFormat = "mm-dd-yy"
$ws.Range("G3").Value = "$ 110"

$ws.Range("A4").Value = "Quiguango Rivera Alisson Lorena"
$ws.Range("B4").Value = 666
$ws.Range("C4").Value = "1003834627"
$ws.Range("D4").Value = "Quiguango Rivera Alisson Lorena"
$ws.Range("E4").Value = "555"
$ws.Range("F4").Value = Get-Date -Year 2020 -Month 11 -Day 2
$ws.Range("F4").NumberFormat = "mm-dd-yy"
$ws.Range("G4").Value = "$ 100"
